$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Text edit: "gave us better results for both LSI and NMF" -> "... for NMF"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "gave us better results for both LSI and NMF",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "gave us better results for NMF", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark away from the end of the document; it will
#    be re-created inside the new hyperlink paragraph below.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3. Insert six new paragraphs at the very top of the document (author name
#    lines, a blank line, a README pointer, a hyperlink line, a blank line).
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
for ($i = 0; $i -lt 6; $i++) {
    $firstPara.Range.InsertParagraphBefore()
}
for ($i = 1; $i -le 6; $i++) {
    $p = $d.Paragraphs($i)
    $p.Style = "Normal"
    $p.Range.Font.Name = "Times New Roman"
}

$d.Paragraphs(1).Range.Text = "Yichen Wu 504294181"
$d.Paragraphs(2).Range.Text = "Siyuan Chen 405024391"
# Paragraph 3 is intentionally left blank.
$d.Paragraphs(4).Range.Text = "Read README.md"
$d.Paragraphs(5).Range.Text = "Or go to "
# Paragraph 6 is intentionally left blank.

# ---------------------------------------------------------------------------
# 4. Build the hyperlink at the end of paragraph 5: "Or go to <link>"
# ---------------------------------------------------------------------------
$linkPara = $d.Paragraphs(5)
$insertPoint = $d.Range($linkPara.Range.End - 1, $linkPara.Range.End - 1)
$hl = $d.Hyperlinks.Add($insertPoint, "https://github.com/popo0293/EE219", $null, $null, "https://github.com/popo0293/EE219")

# Append "Project2" right after the hyperlink text, keeping hyperlink styling.
$hl.Range.InsertAfter("Project2")

# Re-home the "_GoBack" bookmark between the URL and "Project2" (split as
# "P" / "roject2" in the source document), mirroring where the cursor was
# left after the user finished typing the link text.
$markerText = "@@GOBACK_MARK@@"
$hl.Range.InsertAfter($markerText)
$markerFound = $d.Content
$markerFound.Find.Execute($markerText) | Out-Null
$markerStart = $markerFound.Start
$markerRange = $d.Range($markerStart, $markerStart + $markerText.Length)
$markerRange.Text = ""
$bmRange = $d.Range($markerStart, $markerStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
